# Update calibration estimates with values computed using a standard
# moving-average approach (see commit message).

$wb = $excel.ActiveWorkbook

# ---- Linear sheet ----
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.07207301552186338
$wsLinear.Range("B3").Value = 2.760226487825269
$wsLinear.Range("B4").Value = 469.6647525720194
$wsLinear.Range("B5").Value = "[0.9999999999999998, 0.2697471479250945, 0.08056301242197841, 0.10767710410810079, 0.06740160389479591, 0.07291454959977428, 0.2021877883769963, 0.2996385617515471, 0.18863229768646625, 0.03421044222516456, 0.04487798379008036, -0.009599839556339744, -0.003420368383145465, 0.14922228837514626, 0.2761980685595401, 0.1308614903440215, 0.020732076694759256, 0.03080867503577077, 0.018376100044635614, 0.01582687909233444]"

# ---- NonLinear sheet ----
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.458608448726636
$wsNonLinear.Range("B5").Value = 1.723351935365844
$wsNonLinear.Range("B6").Value = 498.6262697683093
$wsNonLinear.Range("B7").Value = -0.7901923414700616
$wsNonLinear.Range("B8").Value = 4.405358692079237
$wsNonLinear.Range("B9").Value = 442.6387877663026
$wsNonLinear.Range("B10").Value = "[1.0, 0.2702649837409783, 0.08150295564028422, 0.10686369009594729, 0.06678014298055358, 0.0733501840236131, 0.20201307761307818, 0.29971925746473255, 0.18902944062853014, 0.034868310929965, 0.04611318814203395, -0.008944689533151173, -0.0022097885189514333, 0.14911228119799902, 0.27538483164479055, 0.13185539188625794, 0.022431954932640216, 0.03207234957281962, 0.017781562171402056, 0.016042924217245715]"
